$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("English", 29.28622855777794),
    @("Japanese", 9.051702747217938),
    @("Spanish", 8.450106595745863),
    @("German", 7.17770091977974),
    @("Chinese", 5.11546168891899),
    @("Arabic", 4.910927738031265),
    @("French", 4.306424068671455),
    @("Portuguese", 4.172165366903836),
    @("Italian", 4.164638326806889),
    @("Malay-Indonesian", 2.373228066497715),
    @("Dutch", 1.850019319643704),
    @("Persian", 1.529928688879911),
    @("Turkish", 1.435392821558656),
    @("Korean", 1.168255134586294),
    @("Polish", 0.9182916681806227),
    @("Thai", 0.8684196098031535),
    @("Urdu", 0.7779000901727845),
    @("Swedish", 0.6205871192406676),
    @("Bengali", 0.404065526279158),
    @("Vietnamese", 0.2807181431933417)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (previously Uzbek/Vietnamese at 22/23)
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(21).Delete()
